$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990207675/1749682509529/zqlwK29jaz-camera-video-segment-driver-1749682514529.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOMZYO4P5%2F20250612%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250612T210013Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJIMEYCIQCw5VXlLF7g5A%2BYIuGL8Een8Dqsh8%2BrEqnjrCzqSmD%2B9QIhAMQDSD00HpcLPuRthGTOtxIlWyTdFsOUS8kVLIwJ8LX%2BKuYDCPX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgxSTL4ZObAP%2BltfwzIqugOm4yr%2B6HVrN7tjK9%2B1FuyRYxFrhaG7RS2OgHJIQkJ0Bz459eqZHtCKy1nhf3CU7%2BUZk07BEDM4%2B5zU8V8N0DfmnKzFgQh0bXTvlKQ0M%2FuUU8mUYONOFoIMMR%2B02bmjurq6zusvJxSLyyRu3OBtt%2FNQHcd49sK5XAt7u8XWPts6%2BUjsdKYe92QPgybVtvKTXOc2EeJsmoTnwCxzwVsGYrna6dYY%2FNVN8yhd3iXHOn1iEbYyEnHOqLqUyZy2PBH%2B6Gn0kqVVUaWtE%2F%2FqFbGvoWCboc5zt9%2F4ZOwApDsbgyqTmT%2BTxZcFmgDaL87YsDxvSapqck%2FBG4biCsgdG4KF5yW2TnV6%2B6B4zg7imvgPVwlyVsyOFuueY%2BoKSs86%2Bi5D4vr4tnqBZ4eAEahBI%2FfAN8WCet16f4Y7s%2BKRlN0lgTJZzbPOb%2F93g8Nbz8qbObXQ0oZBx7Rc3Uz70ZUGGvQBPr%2ByTowh9gfaQgRLu2xqgXKkuM2kN1PN0ldAIutRyTv7aBOnkRpwwzV8xgHwctxZThbI7FRJuDjW1qxC%2BUGgTWHJ1S9M30KRfWPOu90aOeXWSsy3GWls%2Fxh%2BQj0XMObkrMIGOqQBcbr9vE5VAqvTwkUbLtIjEdWH8Bb6bHoEYhL9%2FH%2FRA47U4qRpfEzWtZ33lgF%2BGtP6kLfe%2BkAORjhIWo3Ed%2FSHMdmHimfO7NVfaVrf1JhTmK5cTUoKUPvN4ClZqLI4oDF75%2B3TzzHH3hxsoPY4ol2REUhZEeuB%2FrpD1i%2B05hz7kXDddWOqXIbz4B8%2F1GNGodEg7LAOa%2Bckavl8fhmeTFNMQJCaC0k%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2013%20Jun%202025%2005%3A00%3A13%20GMT&X-Amz-Signature=17f1a16c455c79e888d8c42e59c299dc696b2b1261afa1c0e9f609e8a14a4540'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990207675/1749682509529/JfEAwSUnRy-camera-video-segment-1749682514529.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOMZYO4P5%2F20250612%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250612T210013Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJIMEYCIQCw5VXlLF7g5A%2BYIuGL8Een8Dqsh8%2BrEqnjrCzqSmD%2B9QIhAMQDSD00HpcLPuRthGTOtxIlWyTdFsOUS8kVLIwJ8LX%2BKuYDCPX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgxSTL4ZObAP%2BltfwzIqugOm4yr%2B6HVrN7tjK9%2B1FuyRYxFrhaG7RS2OgHJIQkJ0Bz459eqZHtCKy1nhf3CU7%2BUZk07BEDM4%2B5zU8V8N0DfmnKzFgQh0bXTvlKQ0M%2FuUU8mUYONOFoIMMR%2B02bmjurq6zusvJxSLyyRu3OBtt%2FNQHcd49sK5XAt7u8XWPts6%2BUjsdKYe92QPgybVtvKTXOc2EeJsmoTnwCxzwVsGYrna6dYY%2FNVN8yhd3iXHOn1iEbYyEnHOqLqUyZy2PBH%2B6Gn0kqVVUaWtE%2F%2FqFbGvoWCboc5zt9%2F4ZOwApDsbgyqTmT%2BTxZcFmgDaL87YsDxvSapqck%2FBG4biCsgdG4KF5yW2TnV6%2B6B4zg7imvgPVwlyVsyOFuueY%2BoKSs86%2Bi5D4vr4tnqBZ4eAEahBI%2FfAN8WCet16f4Y7s%2BKRlN0lgTJZzbPOb%2F93g8Nbz8qbObXQ0oZBx7Rc3Uz70ZUGGvQBPr%2ByTowh9gfaQgRLu2xqgXKkuM2kN1PN0ldAIutRyTv7aBOnkRpwwzV8xgHwctxZThbI7FRJuDjW1qxC%2BUGgTWHJ1S9M30KRfWPOu90aOeXWSsy3GWls%2Fxh%2BQj0XMObkrMIGOqQBcbr9vE5VAqvTwkUbLtIjEdWH8Bb6bHoEYhL9%2FH%2FRA47U4qRpfEzWtZ33lgF%2BGtP6kLfe%2BkAORjhIWo3Ed%2FSHMdmHimfO7NVfaVrf1JhTmK5cTUoKUPvN4ClZqLI4oDF75%2B3TzzHH3hxsoPY4ol2REUhZEeuB%2FrpD1i%2B05hz7kXDddWOqXIbz4B8%2F1GNGodEg7LAOa%2Bckavl8fhmeTFNMQJCaC0k%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2013%20Jun%202025%2005%3A00%3A13%20GMT&X-Amz-Signature=f57a4b51cca3e5aa35fa71d5a24f30a2085b5fdf05e1896d361595c53159b078'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991154589/1749681707753/UtIzP4FPML-camera-video-segment-driver-1749681712753.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOMZYO4P5%2F20250612%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250612T210013Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJIMEYCIQCw5VXlLF7g5A%2BYIuGL8Een8Dqsh8%2BrEqnjrCzqSmD%2B9QIhAMQDSD00HpcLPuRthGTOtxIlWyTdFsOUS8kVLIwJ8LX%2BKuYDCPX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgxSTL4ZObAP%2BltfwzIqugOm4yr%2B6HVrN7tjK9%2B1FuyRYxFrhaG7RS2OgHJIQkJ0Bz459eqZHtCKy1nhf3CU7%2BUZk07BEDM4%2B5zU8V8N0DfmnKzFgQh0bXTvlKQ0M%2FuUU8mUYONOFoIMMR%2B02bmjurq6zusvJxSLyyRu3OBtt%2FNQHcd49sK5XAt7u8XWPts6%2BUjsdKYe92QPgybVtvKTXOc2EeJsmoTnwCxzwVsGYrna6dYY%2FNVN8yhd3iXHOn1iEbYyEnHOqLqUyZy2PBH%2B6Gn0kqVVUaWtE%2F%2FqFbGvoWCboc5zt9%2F4ZOwApDsbgyqTmT%2BTxZcFmgDaL87YsDxvSapqck%2FBG4biCsgdG4KF5yW2TnV6%2B6B4zg7imvgPVwlyVsyOFuueY%2BoKSs86%2Bi5D4vr4tnqBZ4eAEahBI%2FfAN8WCet16f4Y7s%2BKRlN0lgTJZzbPOb%2F93g8Nbz8qbObXQ0oZBx7Rc3Uz70ZUGGvQBPr%2ByTowh9gfaQgRLu2xqgXKkuM2kN1PN0ldAIutRyTv7aBOnkRpwwzV8xgHwctxZThbI7FRJuDjW1qxC%2BUGgTWHJ1S9M30KRfWPOu90aOeXWSsy3GWls%2Fxh%2BQj0XMObkrMIGOqQBcbr9vE5VAqvTwkUbLtIjEdWH8Bb6bHoEYhL9%2FH%2FRA47U4qRpfEzWtZ33lgF%2BGtP6kLfe%2BkAORjhIWo3Ed%2FSHMdmHimfO7NVfaVrf1JhTmK5cTUoKUPvN4ClZqLI4oDF75%2B3TzzHH3hxsoPY4ol2REUhZEeuB%2FrpD1i%2B05hz7kXDddWOqXIbz4B8%2F1GNGodEg7LAOa%2Bckavl8fhmeTFNMQJCaC0k%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2013%20Jun%202025%2005%3A00%3A13%20GMT&X-Amz-Signature=b62e71eeb79f21aa3e1ace6187cdab2a38d115e98ea68d3528f4e33131846769'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991154589/1749681707753/418ikFdSWa-camera-video-segment-1749681712753.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOMZYO4P5%2F20250612%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250612T210013Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJIMEYCIQCw5VXlLF7g5A%2BYIuGL8Een8Dqsh8%2BrEqnjrCzqSmD%2B9QIhAMQDSD00HpcLPuRthGTOtxIlWyTdFsOUS8kVLIwJ8LX%2BKuYDCPX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgxSTL4ZObAP%2BltfwzIqugOm4yr%2B6HVrN7tjK9%2B1FuyRYxFrhaG7RS2OgHJIQkJ0Bz459eqZHtCKy1nhf3CU7%2BUZk07BEDM4%2B5zU8V8N0DfmnKzFgQh0bXTvlKQ0M%2FuUU8mUYONOFoIMMR%2B02bmjurq6zusvJxSLyyRu3OBtt%2FNQHcd49sK5XAt7u8XWPts6%2BUjsdKYe92QPgybVtvKTXOc2EeJsmoTnwCxzwVsGYrna6dYY%2FNVN8yhd3iXHOn1iEbYyEnHOqLqUyZy2PBH%2B6Gn0kqVVUaWtE%2F%2FqFbGvoWCboc5zt9%2F4ZOwApDsbgyqTmT%2BTxZcFmgDaL87YsDxvSapqck%2FBG4biCsgdG4KF5yW2TnV6%2B6B4zg7imvgPVwlyVsyOFuueY%2BoKSs86%2Bi5D4vr4tnqBZ4eAEahBI%2FfAN8WCet16f4Y7s%2BKRlN0lgTJZzbPOb%2F93g8Nbz8qbObXQ0oZBx7Rc3Uz70ZUGGvQBPr%2ByTowh9gfaQgRLu2xqgXKkuM2kN1PN0ldAIutRyTv7aBOnkRpwwzV8xgHwctxZThbI7FRJuDjW1qxC%2BUGgTWHJ1S9M30KRfWPOu90aOeXWSsy3GWls%2Fxh%2BQj0XMObkrMIGOqQBcbr9vE5VAqvTwkUbLtIjEdWH8Bb6bHoEYhL9%2FH%2FRA47U4qRpfEzWtZ33lgF%2BGtP6kLfe%2BkAORjhIWo3Ed%2FSHMdmHimfO7NVfaVrf1JhTmK5cTUoKUPvN4ClZqLI4oDF75%2B3TzzHH3hxsoPY4ol2REUhZEeuB%2FrpD1i%2B05hz7kXDddWOqXIbz4B8%2F1GNGodEg7LAOa%2Bckavl8fhmeTFNMQJCaC0k%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2013%20Jun%202025%2005%3A00%3A13%20GMT&X-Amz-Signature=77426c21f91d60966c66696dbd55fc829861d65e01def6bf1d05e7a7450c7ee7'
